$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.481502056121826
$ws.Range("B1").Value = 1.802323460578918
$ws.Range("C1").Value = 1.985083937644958
$ws.Range("D1").Value = 2.256146430969238
$ws.Range("E1").Value = 2.899929761886597
